$wb = $excel.ActiveWorkbook

# Sheet "展览" - F column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6377
$ws1.Range("F3").Value = 111
$ws1.Range("F6").Value = 56
$ws1.Range("F9").Value = 72
$ws1.Range("F12").Value = 155
$ws1.Range("F13").Value = 363
$ws1.Range("F14").Value = 771
$ws1.Range("F15").Value = 3105
$ws1.Range("F18").Value = 1769
$ws1.Range("F19").Value = 21

# Sheet "全部类型" - F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6377
$ws4.Range("F3").Value = 111
$ws4.Range("F6").Value = 56
$ws4.Range("F10").Value = 72
$ws4.Range("F13").Value = 155
$ws4.Range("F14").Value = 363
$ws4.Range("F15").Value = 771
$ws4.Range("F16").Value = 3105
$ws4.Range("F19").Value = 1769
$ws4.Range("F20").Value = 21
